$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.618.01"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.758.57"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.04"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.757.33"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.51%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  -1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.46"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.41%  "

$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("E13").Value = "  -5.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.21"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.387.47"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.753.16"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.568.67"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.02"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.49%  "

$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.00"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.76"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "465.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.697"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.29"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.56%  "

$ws.Range("E25").Value = "  -2.56%  "

$ws.Range("E26").Value = "  -2.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.93"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.05"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.59%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.903.53"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.77"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.35%  "

$ws.Range("E32").Value = "  -4.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.02"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.17"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.16"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.709.66"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.80%  "

$ws.Range("E38").Value = "  -3.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.40"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -8.84%  "

$ws.Range("E40").Value = "  -0.42%  "

$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.78"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.303"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.15"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +8.02%  "

$ws.Range("E47").Value = "  -1.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.49"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.21"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "146.21"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "388.72"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.73%  "
